$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-11-12"

# Update header label in I1 to reflect new "through" date
$ws.Range("I1").Value = "2022 (through 11-12)"

# Update October (row 11), November (row 12) monthly totals and the grand Total (row 14)
$ws.Range("I11").Value = 123
$ws.Range("I12").Value = 35
$ws.Range("I14").Value = 1433
